$wb = $excel.ActiveWorkbook

# === ALC ===
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 300
$ws.Range("J18").Value = 300
$ws.Range("L18").Value = 300
$ws.Range("N18").Value = -868
$ws.Range("H61").Value = 0
$ws.Range("I61").Value = 0
$ws.Range("K61").Value = 0
$ws.Range("M61").ClearContents()
$ws.Range("H64").Value = 19314
$ws.Range("J64").Value = 18856.666
$ws.Range("L64").Value = 18856.666
$ws.Range("N64").Value = -19352.666
$ws.Range("H67").Value = 19314
$ws.Range("J67").Value = 18856.666
$ws.Range("L67").Value = 18856.666
$ws.Range("N67").Value = -20572.666
$ws.Range("H113").Value = 9513.546
$ws.Range("I113").Value = 9404.888999999999
$ws.Range("K113").Value = 9404.888999999999
$ws.Range("M113").Value = -6150.888999999999
$ws.Range("H137").Value = 1452.8334
$ws.Range("I137").Value = 1448.125
$ws.Range("J137").Value = 1462.25
$ws.Range("K137").Value = 4344.375
$ws.Range("L137").Value = 4386.75
$ws.Range("M137").Value = -1794.375
$ws.Range("N137").Value = -9486.75

# === ARM ===
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 900
$ws.Range("I74").Value = 900
$ws.Range("K74").Value = 900
$ws.Range("M74").Value = -26
$ws.Range("H77").Value = 900
$ws.Range("I77").Value = 900
$ws.Range("K77").Value = 4500
$ws.Range("M77").Value = -132
$ws.Range("H88").Value = 3750.5
$ws.Range("I88").Value = 4502.5
$ws.Range("K88").Value = 4502.5
$ws.Range("M88").Value = -4096.5
$ws.Range("H91").Value = 3750.5
$ws.Range("I91").Value = 4502.5
$ws.Range("K91").Value = 4502.5
$ws.Range("M91").Value = -3098.5
$ws.Range("H97").Value = 1109.7142
$ws.Range("I97").Value = 801.4
$ws.Range("K97").Value = 801.4
$ws.Range("M97").Value = -305.4

# === BSM ===
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 1624.4
$ws.Range("I11").Value = 777.3333
$ws.Range("J11").Value = 2895
$ws.Range("K11").Value = 777.3333
$ws.Range("L11").Value = 2895
$ws.Range("M11").Value = -637.3333
$ws.Range("N11").Value = -3175
$ws.Range("H34").Value = 5000
$ws.Range("J34").Value = 5000
$ws.Range("L34").Value = 5000
$ws.Range("N34").Value = -5228
$ws.Range("H86").Value = 4039.6
$ws.Range("I86").Value = 4576
$ws.Range("K86").Value = 4576
$ws.Range("M86").Value = -3453
$ws.Range("H89").Value = 4039.6
$ws.Range("I89").Value = 4576
$ws.Range("K89").Value = 22880
$ws.Range("M89").Value = -17264
$ws.Range("H94").Value = 628.6667
$ws.Range("I94").Value = 628.6667
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 628.6667
$ws.Range("L94").Value = 0
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -177.6667
$ws.Range("H99").Value = 2018.6666
$ws.Range("I99").Value = 1913.9231
$ws.Range("J99").Value = 2699.5
$ws.Range("K99").Value = 1913.9231
$ws.Range("L99").Value = 2699.5
$ws.Range("M99").Value = -415.9231
$ws.Range("N99").Value = -5695.5
$ws.Range("H105").Value = 3337.6365
$ws.Range("I105").Value = 3337.6365
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 3337.6365
$ws.Range("L105").Value = 0
$ws.Range("M105").ClearContents()
$ws.Range("N105").Value = -1590.6365
$ws.Range("H132").Value = 60000
$ws.Range("J132").Value = 60000
$ws.Range("L132").Value = 60000
$ws.Range("N132").Value = -70120
$ws.Range("H134").Value = 1071.5217
$ws.Range("I134").Value = 1106.591
$ws.Range("J134").Value = 300
$ws.Range("K134").Value = 3319.773
$ws.Range("L134").Value = 900
$ws.Range("M134").Value = -784.7729999999997
$ws.Range("N134").Value = -5970

# === CRP ===
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 3957.875
$ws.Range("I31").Value = 3410.5
$ws.Range("K31").Value = 3410.5
$ws.Range("M31").Value = -3115.5
$ws.Range("H34").Value = 3957.875
$ws.Range("I34").Value = 3410.5
$ws.Range("K34").Value = 3410.5
$ws.Range("M34").Value = -3208.5
$ws.Range("H43").Value = 23986.666
$ws.Range("J43").Value = 23986.666
$ws.Range("L43").Value = 23986.666
$ws.Range("N43").Value = -24354.666
$ws.Range("H52").Value = 35000
$ws.Range("I52").Value = 0
$ws.Range("J52").Value = 35000
$ws.Range("K52").Value = 0
$ws.Range("L52").ClearContents()
$ws.Range("M52").Value = 35000
$ws.Range("N52").Value = -35588
$ws.Range("H95").Value = 13263.857
$ws.Range("J95").Value = 13263.857
$ws.Range("L95").Value = 13263.857
$ws.Range("N95").Value = -18755.857
$ws.Range("H101").Value = 23986.666
$ws.Range("J101").Value = 23986.666
$ws.Range("L101").Value = 23986.666
$ws.Range("N101").Value = -30476.666
$ws.Range("H134").Value = 1792.3636
$ws.Range("I134").Value = 1335.7778
$ws.Range("K134").Value = 4007.3334
$ws.Range("M134").Value = -1472.3334

# === CUL ===
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 392.5
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 392.5
$ws.Range("K23").Value = 0
$ws.Range("L23").ClearContents()
$ws.Range("M23").Value = 1177.5
$ws.Range("N23").Value = -1647.5
$ws.Range("H113").Value = 1180.4445
$ws.Range("I113").Value = 603.3333
$ws.Range("J113").Value = 1469
$ws.Range("K113").Value = 1809.9999
$ws.Range("L113").Value = 4407
$ws.Range("M113").Value = 360.0001
$ws.Range("N113").Value = -8747
$ws.Range("H131").Value = 1001.7143
$ws.Range("H140").Value = 569.44446
$ws.Range("I140").Value = 569.44446
$ws.Range("K140").Value = 1708.33338
$ws.Range("M140").Value = 3471.66662

# === GSM ===
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H25").Value = 601
$ws.Range("J25").Value = 601
$ws.Range("L25").Value = 601
$ws.Range("N25").Value = -1659
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").ClearContents()
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = 0
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").ClearContents()
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = 0
$ws.Range("H102").Value = 2624.125
$ws.Range("I102").Value = 2417.1667
$ws.Range("J102").Value = 3245
$ws.Range("K102").Value = 2417.1667
$ws.Range("L102").Value = 3245
$ws.Range("M102").Value = -795.1667000000002
$ws.Range("N102").Value = -6489
$ws.Range("H107").Value = 3534.7646
$ws.Range("I107").Value = 501.83334
$ws.Range("J107").Value = 10813.8
$ws.Range("K107").Value = 501.83334
$ws.Range("L107").Value = 10813.8
$ws.Range("M107").Value = 1418.16666
$ws.Range("N107").Value = -14653.8
$ws.Range("H122").Value = 0
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 0
$ws.Range("L122").ClearContents()
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = 0
$ws.Range("H132").Value = 2843.4546
$ws.Range("I132").Value = 2886.4
$ws.Range("K132").Value = 8659.200000000001
$ws.Range("M132").Value = -6129.200000000001
$ws.Range("H137").Value = 79443
$ws.Range("J137").Value = 79443
$ws.Range("L137").Value = 79443
$ws.Range("N137").Value = -89643

# === LTW ===
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 1997.5
$ws.Range("I7").Value = 1997.5
$ws.Range("K7").Value = 1997.5
$ws.Range("M7").Value = -1885.5
$ws.Range("H22").Value = 1866.6666
$ws.Range("H27").Value = 1866.6666
$ws.Range("H40").Value = 2500
$ws.Range("I40").Value = 2500
$ws.Range("K40").Value = 2500
$ws.Range("M40").Value = -2364
$ws.Range("H46").Value = 4376.963
$ws.Range("J46").Value = 4679.933
$ws.Range("L46").Value = 4679.933
$ws.Range("N46").Value = -5055.933
$ws.Range("H68").Value = 3875
$ws.Range("J68").Value = 4750
$ws.Range("L68").Value = 4750
$ws.Range("N68").Value = -6248
$ws.Range("H71").Value = 3875
$ws.Range("J71").Value = 4750
$ws.Range("L71").Value = 23750
$ws.Range("N71").Value = -31238
$ws.Range("H95").Value = 9999.5
$ws.Range("J95").Value = 9999.5
$ws.Range("L95").Value = 9999.5
$ws.Range("N95").Value = -15491.5
$ws.Range("H126").Value = 1997.5
$ws.Range("I126").Value = 1997.5
$ws.Range("K126").Value = 5992.5
$ws.Range("M126").Value = -3522.5
$ws.Range("H136").Value = 1768.238
$ws.Range("I136").Value = 1815.7
$ws.Range("K136").Value = 5447.1
$ws.Range("M136").Value = -2897.1
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").ClearContents()
$ws.Range("N137").Value = 0

# === WVR ===
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H58").Value = 0
$ws.Range("I58").Value = 0
$ws.Range("K58").Value = 0
$ws.Range("M58").ClearContents()
$ws.Range("H136").Value = 1102.5
$ws.Range("I136").Value = 1163.8462
$ws.Range("K136").Value = 3491.5386
$ws.Range("M136").Value = -941.5385999999999
